$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.354.88'
$ws.Range("E2").Value = '  +2.97%  '
$ws.Range("D3").Value = '1.724.24'
$ws.Range("E3").Value = '  +3.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.37'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4718'
$ws.Range("E7").Value = '  -2.20%  '
$ws.Range("E8").Value = '  -0.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06226'
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("D10").Value = '1.718.91'
$ws.Range("E10").Value = '  +3.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07092'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.35'
$ws.Range("E12").Value = '  +2.58%  '
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.402'
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.37'
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = '26.352.49'
$ws.Range("E18").Value = '  +3.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006795'
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.59'
$ws.Range("E20").Value = '  +0.97%  '
$ws.Range("D21").Value = '1.938.69'
$ws.Range("E21").Value = '  +3.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.561'
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.765'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.345'
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.00'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.25'
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '108.75'
$ws.Range("E27").Value = '  +3.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.409'
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.767'
$ws.Range("E29").Value = '  +3.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.013'
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.699'
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07756'
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04467'
$ws.Range("E33").Value = '  +2.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.613'
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9791'
$ws.Range("E35").Value = '  +2.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6209'
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '115.40'
$ws.Range("E37").Value = '  +17.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9223'
$ws.Range("E38").Value = '  +6.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.426'
$ws.Range("E39").Value = '  -7.50%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.000'
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.906'
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01479'
$ws.Range("E42").Value = '  -2.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.391'
$ws.Range("E43").Value = '  +14.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3822'
$ws.Range("E44").Value = '  +0.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1166'
$ws.Range("E45").Value = '  +3.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.265'
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("E47").Value = '  +0.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.67'
$ws.Range("E48").Value = '  +3.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.693'
$ws.Range("E49").Value = '  +3.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3392'
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.220'
$ws.Range("E51").Value = '  +1.54%  '
Write-Output "Applied cryptos update"
